# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with
# the latest scraped values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal TEXT (not an auto-converted
# number/date), matching the inline-string cells already in the sheet,
# while preserving that cell's original style/formatting.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "61.914.88"
$ws.Range("E2").Value = "  -1.58%  "

Set-TextValue $ws.Range("D3") "3.014.94"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue $ws.Range("D5") "542.29"
$ws.Range("E5").Value = "  +1.03%  "

Set-TextValue $ws.Range("D6") "133.76"
$ws.Range("E6").Value = "  -2.14%  "

$ws.Range("E7").Value = "  +0.07%  "

Set-TextValue $ws.Range("D8") "3.008.42"
$ws.Range("E8").Value = "  -1.51%  "

$ws.Range("E9").Value = "  -0.08%  "

Set-TextValue $ws.Range("D10") "6.22"
$ws.Range("E10").Value = "  -0.23%  "

Set-TextValue $ws.Range("D11") "0.147"
$ws.Range("E11").Value = "  -5.24%  "

$ws.Range("E12").Value = "  -1.29%  "

Set-TextValue $ws.Range("D13") "34.63"
$ws.Range("E13").Value = "  +0.79%  "

Set-TextValue $ws.Range("D14") "0.0000220"
$ws.Range("E14").Value = "  -0.54%  "

Set-TextValue $ws.Range("D15") "3.500.72"
$ws.Range("E15").Value = "  -1.53%  "

Set-TextValue $ws.Range("D16") "61.895.31"
$ws.Range("E16").Value = "  -1.59%  "

Set-TextValue $ws.Range("D17") "0.110"
$ws.Range("E17").Value = "  -2.60%  "

Set-TextValue $ws.Range("D18") "3.005.12"
$ws.Range("E18").Value = "  -1.70%  "

Set-TextValue $ws.Range("D19") "6.63"
$ws.Range("E19").Value = "  +0.47%  "

Set-TextValue $ws.Range("D20") "482.81"
$ws.Range("E20").Value = "  +3.12%  "

Set-TextValue $ws.Range("D21") "13.25"
$ws.Range("E21").Value = "  -1.61%  "

Set-TextValue $ws.Range("D22") "0.671"
$ws.Range("E22").Value = "  -3.23%  "

Set-TextValue $ws.Range("D23") "6.98"
$ws.Range("E23").Value = "  -0.20%  "

Set-TextValue $ws.Range("D24") "81.92"
$ws.Range("E24").Value = "  +4.69%  "

Set-TextValue $ws.Range("D25") "11.98"
$ws.Range("E25").Value = "  -0.61%  "

$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("E27").Value = "  +0.42%  "

$ws.Range("E28").Value = "  -1.72%  "

Set-TextValue $ws.Range("D29") "0.997"
$ws.Range("E29").Value = "  -0.24%  "

Set-TextValue $ws.Range("D30") "1.92"
$ws.Range("E30").Value = "  +2.64%  "

Set-TextValue $ws.Range("D31") "25.71"
$ws.Range("E31").Value = "  -1.19%  "

Set-TextValue $ws.Range("D32") "1.13"
$ws.Range("E32").Value = "  -1.14%  "

Set-TextValue $ws.Range("D33") "5.67"
$ws.Range("E33").Value = "  +3.96%  "

Set-TextValue $ws.Range("D34") "2.34"
$ws.Range("E34").Value = "  +1.82%  "

Set-TextValue $ws.Range("D35") "55.31"
$ws.Range("E35").Value = "  -6.23%  "

Set-TextValue $ws.Range("D36") "5.87"
$ws.Range("E36").Value = "  -1.20%  "

Set-TextValue $ws.Range("D37") "446.20"
$ws.Range("E37").Value = "  -7.21%  "

Set-TextValue $ws.Range("D38") "3.150.12"
$ws.Range("E38").Value = "  -2.77%  "

Set-TextValue $ws.Range("D39") "0.0796"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("E40").Value = "  -3.37%  "

$ws.Range("E41").Value = "  +0.51%  "

Set-TextValue $ws.Range("D42") "8.08"
$ws.Range("E42").Value = "  -0.25%  "

Set-TextValue $ws.Range("D43") "2.44"
$ws.Range("E43").Value = "  -4.48%  "

Set-TextValue $ws.Range("D44") "26.38"
$ws.Range("E44").Value = "  +5.29%  "

Set-TextValue $ws.Range("D46") "0.242"
$ws.Range("E46").Value = "  -3.10%  "

$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("E48").Value = "  -2.01%  "

Set-TextValue $ws.Range("D49") "115.92"
$ws.Range("E49").Value = "  -5.88%  "

$ws.Range("E50").Value = "  +5.14%  "

Set-TextValue $ws.Range("D51") "0.0₃0489"
$ws.Range("E51").Value = "  -5.49%  "
